$d = $word.ActiveDocument

# The target paragraph is the 3rd paragraph in the document body: an
# otherwise-empty paragraph (paraId 6C7F2931 in the original OOXML) that
# sits between the title block and the "Objetivo:" paragraph. We need to
# add a new run containing two spaces, formatted with the same run
# properties already carried on the paragraph mark (rFonts eastAsiaTheme,
# a gray theme color, and size 24).
$para = $d.Paragraphs(3)
$rng = $para.Range

# Use a zero-length range positioned at the very start of the paragraph
# (strictly inside its content, not at the paragraph-mark boundary) so
# InsertXML merges the new run into the existing paragraph instead of
# replacing the whole <w:p> (which would otherwise drop its w:pPr/paraId).
$insertionPoint = $d.Range($rng.Start, $rng.Start)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/part.xml" pkg:contentType="application/xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p><w:r>' + `
  '<w:rPr>' + `
  '<w:rFonts w:eastAsiaTheme="majorEastAsia"/>' + `
  '<w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/>' + `
  '<w:sz w:val="24"/>' + `
  '<w:szCs w:val="24"/>' + `
  '</w:rPr>' + `
  '<w:t xml:space="preserve">  </w:t>' + `
  '</w:r></w:p></w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml)
